$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("pcodeFile")

# Insert a new column B ("qa code") in front of the existing "method code" column,
# shifting all the other pcode columns one place to the right.
$ws1.Columns.Item(2).Insert() | Out-Null
$ws1.Cells.Item(1, 2).Value = "qa code"

# Give the new column the same width as column A (matches the exported
# width of 6.5703125 used by column A).
$ws1.Columns.Item(2).ColumnWidth = $ws1.Columns.Item(1).ColumnWidth

# Make "pcodeFile" the active sheet/tab, with B2 selected, instead of
# "qwsampleFile" being active.
$ws1.Activate() | Out-Null
$ws1.Range("B2").Select() | Out-Null
